$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.879.21"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "2.667.72"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.650"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.26%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.126"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.402"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000194"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").Value = "3.148.83"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "65.788.61"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").Value = "2.704.30"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.46%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("E24").Value = "  +11.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000112"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "566.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.163"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.27%  "
$ws.Range("E35").Value = "  +0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.423"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "161.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0619"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.646"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.64%  "
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("E50").Value = "  -6.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.814"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.82%  "
